# Weekly update: insert 3 new "Haba" price records into the Feria Lagunitas
# de Puerto Montt sheet, shifting existing rows down as needed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-HabaRow {
    param([int]$Row, [double]$Fecha, [double]$Volumen, [double]$PrecioMinimo, [double]$PrecioMaximo, [double]$PrecioPromedio, [string]$Origen, [double]$PrecioKg)

    $ws.Cells.Item($Row, 1).Value = 4
    $ws.Cells.Item($Row, 2).Value = "Feria Lagunitas de Puerto Montt"
    $ws.Cells.Item($Row, 3).Value = "Los Lagos"
    $ws.Cells.Item($Row, 4).Value = $Fecha
    $ws.Cells.Item($Row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($Row, 5).Value = 10
    $ws.Cells.Item($Row, 6).Value = 100112026
    $ws.Cells.Item($Row, 7).Value = "Haba"
    $ws.Cells.Item($Row, 8).Value = "Sin especificar"
    $ws.Cells.Item($Row, 9).Value = "Primera"
    $ws.Cells.Item($Row, 10).Value = $Volumen
    $ws.Cells.Item($Row, 11).Value = $PrecioMinimo
    $ws.Cells.Item($Row, 12).Value = $PrecioMaximo
    $ws.Cells.Item($Row, 13).Value = $PrecioPromedio
    $ws.Cells.Item($Row, 14).Value = "$/saco 25 kilos"
    $ws.Cells.Item($Row, 15).Value = $Origen
    $ws.Cells.Item($Row, 16).Value = $PrecioKg
    $ws.Cells.Item($Row, 17).Value = 25
    $ws.Cells.Item($Row, 18).Value = "Hortaliza"
}

# 1) Insert a new record (fecha 44425) above the row currently holding 44392.
$ws.Rows.Item(5).Insert()
Set-HabaRow 5 44425 90 18000 18000 18000 "Provincia de Limarí" 720

# 2) Insert a new record (fecha 44421) above the row currently holding 44400.
$ws.Rows.Item(17).Insert()
Set-HabaRow 17 44421 80 17000 17000 17000 "Provincia de Limarí" 680

# 3) Append a new record (fecha 44418) at the end of the table.
Set-HabaRow 31 44418 90 18000 18000 18000 "Provincia de Limarí" 720
